# Fruta / hortaliza, semanal
#
# The sheet contains one new weekly observation that needs to be inserted
# as row 492 ("Dulce o Americano" / "Primera", fecha serial 44694 =
# 2022-05-13). Inserting the row pushes every existing row from 492
# through 593 down by one (to 493 through 594), and Excel's native
# row-insert behaviour takes care of re-numbering / re-basing all of the
# existing data automatically, so all we need to do is insert the row and
# then populate the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 492; rows 492:593 shift down to 493:594
$ws.Rows(492).Insert()

# Populate the newly inserted row with the new observation
$ws.Range("A492").Value = 8
$ws.Range("B492").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C492").Value = 'Coquimbo'
$ws.Range("D492").Value = 44694
$ws.Range("E492").Value = 4
$ws.Range("F492").Value = 100112024
$ws.Range("G492").Value = 'Choclo'
$ws.Range("H492").Value = 'Dulce o Americano'
$ws.Range("I492").Value = 'Primera'
$ws.Range("J492").Value = 22000
$ws.Range("K492").Value = 300
$ws.Range("L492").Value = 350
$ws.Range("M492").Value = 325
$ws.Range("N492").Value = '$/unidad'
$ws.Range("O492").Value = 'Provincia del Elquí'
$ws.Range("P492").Value = 325
$ws.Range("Q492").Value = 1
$ws.Range("R492").Value = 'Hortaliza'
